# FALECPV-CuentaCobrar.xlsx - "Cuentas por Cobrar" sheet fix
#
# The "Cuentas por Cobrar" (accounts receivable) header row had a column
# labelled ADEUDADO (owed) immediately followed by ABONO (payment) and then
# SALDO (balance) / VENCIMIENTO (due date). The fix relabels things so the
# columns read: ... CRÉDITO, ABONO, ADEUDADO, VENCIMIENTO - i.e. the old
# "ADEUDADO" header becomes "CRÉDITO" and the old "SALDO" header becomes
# "ADEUDADO" (the real running total column), fixing the totals row per the
# commit message ("pantalla de cuentas por cobrar arreglar totales").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 9): I9 was "ADEUDADO" -> now "CRÉDITO"
$ws.Range("I9").Value = "CRÉDITO"

# Header row (row 9): K9 was "SALDO" -> now "ADEUDADO"
$ws.Range("K9").Value = "ADEUDADO"

# Reflect the last-active cell from the authoring session.
$ws.Range("K9").Select()
